# Update PLC live data values on the active sheet ("LiveData")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 285
$ws.Range("C3").Value = 175315
$ws.Range("C4").Value = 165288
$ws.Range("C5").Value = 10028
$ws.Range("C8").Value = 64.44
